$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos text replaced by docente name ---
$ws.Range("B10").Value = '8711686 - Flavia Reis Cardoso Rojas'
$ws.Range("C10").Value = '8711686 - Flavia Reis Cardoso Rojas'

# --- Row 13: new A label + Semestral replacing the docente line ---
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

# --- Row 14: relabeled to Short syllabus: (B/C keep their text) ---
$ws.Range("A14").Value = 'Short syllabus:'

# --- Row 15: relabeled to Programa:, value becomes the activation date string ---
$ws.Range("A15").Value = 'Programa:'
# Copy/PasteSpecial(values) from B8/C8 (already text-typed "01/01/2018") so the
# cell stays a shared-string text cell instead of being reinterpreted as a date.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Row 16: relabeled Syllabus:, English syllabus text ---
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton''s laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum.'
$ws.Range("C16").Value = '1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton''s laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum.'

# --- Row 17: relabeled Avaliação:, B/C no longer used ---
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# --- Row 18: relabeled Método:, new B/C with docente name ---
# B18/C18 did not exist before, so give them the same cell format as the
# other B/C column cells (style index 2 / 3) before writing the text -
# otherwise a plain Value= on a previously-blank cell falls back to a
# bold/no-wrap style and a brand new style entry gets created.
$ws.Range("A18").Value = 'Método:'
$ws.Range("B13").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B18").Value = '8711686 - Flavia Reis Cardoso Rojas'
$ws.Range("C18").Value = '8711686 - Flavia Reis Cardoso Rojas'

# --- Rows 19-21: each label shifts to the next item (B/C text unchanged) ---
$ws.Range("A19").Value = 'Critério:'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("A21").Value = 'Bibliografia:'

# --- Remove the old trailing row (detailed bibliography entry) ---
$ws.Rows.Item(22).Delete()

# --- Row height adjustments to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
